$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1985559566787004
$ws.Range("C2").Value = 0.5415162454873647
$ws.Range("J2").Value = 0.02888086642599278
$ws.Range("P2").Value = 0.1407942238267148
$ws.Range("S2").Value = 0.09025270758122744
$ws.Range("B3").Value = 0.00641025641025641
$ws.Range("C3").Value = 0.01923076923076923
$ws.Range("J3").Value = 0.02564102564102564
$ws.Range("O3").Value = 0.00641025641025641
$ws.Range("P3").Value = 0.7564102564102564
$ws.Range("S3").Value = 0.1858974358974359
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.03488372093023256
$ws.Range("D6").Value = 0.01162790697674419
$ws.Range("F6").Value = 0.04651162790697674
$ws.Range("J6").Value = 0.2267441860465116
$ws.Range("O6").Value = 0.01744186046511628
$ws.Range("Q6").Value = 0.1395348837209302
$ws.Range("R6").Value = 0.1220930232558139
$ws.Range("S6").Value = 0.4011627906976744
$ws.Range("B7").Value = 0.1694915254237288
$ws.Range("D7").Value = 0.01694915254237288
$ws.Range("F7").Value = 0.05084745762711865
$ws.Range("J7").Value = 0.1694915254237288
$ws.Range("O7").Value = 0.008474576271186441
$ws.Range("Q7").Value = 0.1694915254237288
$ws.Range("R7").Value = 0.0847457627118644
$ws.Range("S7").Value = 0.3305084745762712
$ws.Range("B8").Value = 0.09798270893371758
$ws.Range("D8").Value = 0.01440922190201729
$ws.Range("F8").Value = 0.05475504322766571
$ws.Range("J8").Value = 0.09510086455331412
$ws.Range("O8").Value = 0.01152737752161383
$ws.Range("Q8").Value = 0.207492795389049
$ws.Range("R8").Value = 0.138328530259366
$ws.Range("S8").Value = 0.3804034582132565
$ws.Range("B9").Value = 0.1878787878787879
$ws.Range("D9").Value = 0.006060606060606061
$ws.Range("F9").Value = 0.05454545454545454
$ws.Range("J9").Value = 0.09090909090909091
$ws.Range("O9").Value = 0.01818181818181818
$ws.Range("Q9").Value = 0.1575757575757576
$ws.Range("R9").Value = 0.08484848484848485
$ws.Range("S9").Value = 0.4
$ws.Range("B10").Value = 0.1248799231508165
$ws.Range("D10").Value = 0.02401536983669549
$ws.Range("E10").Value = 0.0009606147934678194
$ws.Range("F10").Value = 0.07684918347742556
$ws.Range("J10").Value = 0.1018251681075889
$ws.Range("O10").Value = 0.009606147934678195
$ws.Range("Q10").Value = 0.1882804995196926
$ws.Range("R10").Value = 0.1268011527377522
$ws.Range("S10").Value = 0.3467819404418828
$ws.Range("G11").Value = 0.1929824561403509
$ws.Range("J11").Value = 0.06432748538011696
$ws.Range("K11").Value = 0.2222222222222222
$ws.Range("L11").Value = 0.4970760233918128
$ws.Range("S11").Value = 0.02339181286549707
$ws.Range("G12").Value = 0.7078651685393258
$ws.Range("J12").Value = 0.2022471910112359
$ws.Range("K12").Value = 0.02247191011235955
$ws.Range("L12").Value = 0.03370786516853932
$ws.Range("S12").Value = 0.03370786516853932
$ws.Range("G13").Value = 0.6410256410256411
$ws.Range("J13").Value = 0.2564102564102564
$ws.Range("S13").Value = 0.1025641025641026
$ws.Range("F15").Value = 0.005847953216374269
$ws.Range("H15").Value = 0.1695906432748538
$ws.Range("I15").Value = 0.1169590643274854
$ws.Range("J15").Value = 0.3742690058479532
$ws.Range("K15").Value = 0.01754385964912281
$ws.Range("M15").Value = 0.01169590643274854
$ws.Range("N15").Value = 0.01169590643274854
$ws.Range("O15").Value = 0.06432748538011696
$ws.Range("S15").Value = 0.2280701754385965
$ws.Range("F16").Value = 0.005681818181818182
$ws.Range("H16").Value = 0.1761363636363636
$ws.Range("I16").Value = 0.07386363636363637
$ws.Range("J16").Value = 0.3636363636363636
$ws.Range("K16").Value = 0.08522727272727272
$ws.Range("M16").Value = 0.02840909090909091
$ws.Range("O16").Value = 0.1079545454545455
$ws.Range("S16").Value = 0.1590909090909091
$ws.Range("F17").Value = 0.01515151515151515
$ws.Range("H17").Value = 0.1484848484848485
$ws.Range("I17").Value = 0.1272727272727273
$ws.Range("J17").Value = 0.4151515151515152
$ws.Range("K17").Value = 0.07272727272727272
$ws.Range("M17").Value = 0.01818181818181818
$ws.Range("O17").Value = 0.0696969696969697
$ws.Range("S17").Value = 0.1333333333333333
$ws.Range("F18").Value = 0.03603603603603604
$ws.Range("H18").Value = 0.1801801801801802
$ws.Range("I18").Value = 0.06756756756756757
$ws.Range("J18").Value = 0.4414414414414414
$ws.Range("K18").Value = 0.06756756756756757
$ws.Range("M18").Value = 0.02252252252252252
$ws.Range("O18").Value = 0.06756756756756757
$ws.Range("S18").Value = 0.1171171171171171
$ws.Range("F19").Value = 0.009027081243731194
$ws.Range("H19").Value = 0.1995987963891675
$ws.Range("I19").Value = 0.07622868605817452
$ws.Range("J19").Value = 0.4252758274824474
$ws.Range("K19").Value = 0.07321965897693079
$ws.Range("M19").Value = 0.02407221664994985
$ws.Range("O19").Value = 0.0641925777331996
$ws.Range("S19").Value = 0.1283851554663992
